$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to update: row -> B, C, D, E values
$data = @{
    2  = @("30.23", "39.36", "2024-10", "2024-10")
    3  = @("22.11", "43.82", "2024-10", "2024-10")
    4  = @("40.61", "50.41", "2024-10", "2024-10")
    5  = @("27.07", "36.63", "2024-10", "2024-10")
    6  = @("37.86", "39.03", "2024-10", "2024-10")
    7  = @("29.54", "37.34", "2024-10", "2024-10")
    8  = @("37.95", "37.91", "2024-10", "2024-10")
    9  = @("22.22", "41.73", "2024-10", "2024-10")
    10 = @("27.17", "39.42", "2024-10", "2024-10")
    11 = @("35.98", "40.94", "2024-10", "2024-10")
}

# Force text format on the whole data range first so numeric-looking
# strings (e.g. "30.23") stay as text instead of being parsed as numbers.
$ws.Range("B2:E11").NumberFormat = "@"

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
}
